$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "75÷9="
$t.Cell(1, 2).Range.Text = "50÷5="
$t.Cell(1, 3).Range.Text = "14÷3="
$t.Cell(1, 4).Range.Text = "29÷2="
$t.Cell(1, 5).Range.Text = "59÷7="

# Row 5
$t.Cell(5, 1).Range.Text = "32÷6="
$t.Cell(5, 2).Range.Text = "33÷8="
$t.Cell(5, 3).Range.Text = "34÷5="
$t.Cell(5, 4).Range.Text = "80÷5="
$t.Cell(5, 5).Range.Text = "25÷5="

# Row 9
$t.Cell(9, 1).Range.Text = "53÷5="
$t.Cell(9, 2).Range.Text = "46÷5="
$t.Cell(9, 3).Range.Text = "89÷7="
$t.Cell(9, 4).Range.Text = "75÷9="
$t.Cell(9, 5).Range.Text = "83÷3="

# Row 13
$t.Cell(13, 1).Range.Text = "58÷5="
$t.Cell(13, 2).Range.Text = "51÷5="
$t.Cell(13, 3).Range.Text = "94÷2="
$t.Cell(13, 4).Range.Text = "95÷8="
$t.Cell(13, 5).Range.Text = "54÷2="

# Row 17
$t.Cell(17, 1).Range.Text = "94÷4="
$t.Cell(17, 2).Range.Text = "60÷6="
$t.Cell(17, 3).Range.Text = "88÷8="
$t.Cell(17, 4).Range.Text = "12÷5="
$t.Cell(17, 5).Range.Text = "43÷8="
